$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-05-20 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-05-21 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("72-6=66", $true, $false, $false, $false, $false, $true, 1, $false, "29-2=27", 2) | Out-Null
$d.Content.Find.Execute("26+41=67", $true, $false, $false, $false, $false, $true, 1, $false, "2+66=68", 2) | Out-Null
$d.Content.Find.Execute("50+29=79", $true, $false, $false, $false, $false, $true, 1, $false, "46-1=45", 2) | Out-Null
$d.Content.Find.Execute("27+11=38", $true, $false, $false, $false, $false, $true, 1, $false, "93-40=53", 2) | Out-Null
$d.Content.Find.Execute("75+3=78", $true, $false, $false, $false, $false, $true, 1, $false, "61-44=17", 2) | Out-Null
$d.Content.Find.Execute("83-36=47", $true, $false, $false, $false, $false, $true, 1, $false, "80-20=60", 2) | Out-Null
$d.Content.Find.Execute("84-55=29", $true, $false, $false, $false, $false, $true, 1, $false, "38+27=65", 2) | Out-Null
$d.Content.Find.Execute("73-10=63", $true, $false, $false, $false, $false, $true, 1, $false, "52-17=35", 2) | Out-Null
$d.Content.Find.Execute("3+72=75", $true, $false, $false, $false, $false, $true, 1, $false, "72+5=77", 2) | Out-Null
$d.Content.Find.Execute("66+8=74", $true, $false, $false, $false, $false, $true, 1, $false, "8+60=68", 2) | Out-Null
$d.Content.Find.Execute("44-5=39", $true, $false, $false, $false, $false, $true, 1, $false, "38+47=85", 2) | Out-Null
$d.Content.Find.Execute("10+44=54", $true, $false, $false, $false, $false, $true, 1, $false, "27+5=32", 2) | Out-Null
$d.Content.Find.Execute("4+54=58", $true, $false, $false, $false, $false, $true, 1, $false, "39-4=35", 2) | Out-Null
$d.Content.Find.Execute("41+8=49", $true, $false, $false, $false, $false, $true, 1, $false, "17+59=76", 2) | Out-Null
$d.Content.Find.Execute("4+42=46", $true, $false, $false, $false, $false, $true, 1, $false, "76-11=65", 2) | Out-Null
$d.Content.Find.Execute("19-12=7", $true, $false, $false, $false, $false, $true, 1, $false, "50-23=27", 2) | Out-Null
$d.Content.Find.Execute("22+45=67", $true, $false, $false, $false, $false, $true, 1, $false, "48+51=99", 2) | Out-Null
$d.Content.Find.Execute("77+14=91", $true, $false, $false, $false, $false, $true, 1, $false, "38-12=26", 2) | Out-Null
$d.Content.Find.Execute("21-7=14", $true, $false, $false, $false, $false, $true, 1, $false, "66+20=86", 2) | Out-Null
$d.Content.Find.Execute("72-37=35", $true, $false, $false, $false, $false, $true, 1, $false, "74+18=92", 2) | Out-Null
$d.Content.Find.Execute("37+4=41", $true, $false, $false, $false, $false, $true, 1, $false, "11+2=13", 2) | Out-Null
$d.Content.Find.Execute("46+50=96", $true, $false, $false, $false, $false, $true, 1, $false, "60-57=3", 2) | Out-Null
$d.Content.Find.Execute("17+16=33", $true, $false, $false, $false, $false, $true, 1, $false, "83-24=59", 2) | Out-Null
$d.Content.Find.Execute("33-0=33", $true, $false, $false, $false, $false, $true, 1, $false, "30+16=46", 2) | Out-Null
$d.Content.Find.Execute("78-70=8", $true, $false, $false, $false, $false, $true, 1, $false, "10+33=43", 2) | Out-Null
$d.Content.Find.Execute("71-28=43", $true, $false, $false, $false, $false, $true, 1, $false, "42+15=57", 2) | Out-Null
$d.Content.Find.Execute("23+1=24", $true, $false, $false, $false, $false, $true, 1, $false, "13+10=23", 2) | Out-Null
$d.Content.Find.Execute("71+3=74", $true, $false, $false, $false, $false, $true, 1, $false, "59+33=92", 2) | Out-Null
$d.Content.Find.Execute("68-55=13", $true, $false, $false, $false, $false, $true, 1, $false, "75-74=1", 2) | Out-Null
$d.Content.Find.Execute("28+64=92", $true, $false, $false, $false, $false, $true, 1, $false, "61-42=19", 2) | Out-Null
$d.Content.Find.Execute("31+52=83", $true, $false, $false, $false, $false, $true, 1, $false, "19-7=12", 2) | Out-Null
$d.Content.Find.Execute("91-38=53", $true, $false, $false, $false, $false, $true, 1, $false, "23+68=91", 2) | Out-Null
$d.Content.Find.Execute("4+82=86", $true, $false, $false, $false, $false, $true, 1, $false, "60+33=93", 2) | Out-Null
$d.Content.Find.Execute("25-12=13", $true, $false, $false, $false, $false, $true, 1, $false, "53-17=36", 2) | Out-Null
$d.Content.Find.Execute("3+42=45", $true, $false, $false, $false, $false, $true, 1, $false, "79-73=6", 2) | Out-Null
$d.Content.Find.Execute("3+68=71", $true, $false, $false, $false, $false, $true, 1, $false, "99-17=82", 2) | Out-Null
$d.Content.Find.Execute("42-36=6", $true, $false, $false, $false, $false, $true, 1, $false, "97-58=39", 2) | Out-Null
$d.Content.Find.Execute("39-19=20", $true, $false, $false, $false, $false, $true, 1, $false, "88-56=32", 2) | Out-Null
$d.Content.Find.Execute("84-34=50", $true, $false, $false, $false, $false, $true, 1, $false, "94-49=45", 2) | Out-Null
$d.Content.Find.Execute("55-43=12", $true, $false, $false, $false, $false, $true, 1, $false, "33+36=69", 2) | Out-Null
$d.Content.Find.Execute("55+16=71", $true, $false, $false, $false, $false, $true, 1, $false, "67+4=71", 2) | Out-Null
$d.Content.Find.Execute("36-28=8", $true, $false, $false, $false, $false, $true, 1, $false, "62-21=41", 2) | Out-Null
$d.Content.Find.Execute("84-82=2", $true, $false, $false, $false, $false, $true, 1, $false, "51+16=67", 2) | Out-Null
$d.Content.Find.Execute("59+3=62", $true, $false, $false, $false, $false, $true, 1, $false, "97-91=6", 2) | Out-Null
$d.Content.Find.Execute("37+58=95", $true, $false, $false, $false, $false, $true, 1, $false, "0+61=61", 2) | Out-Null
$d.Content.Find.Execute("54-10=44", $true, $false, $false, $false, $false, $true, 1, $false, "94-28=66", 2) | Out-Null
$d.Content.Find.Execute("22-19=3", $true, $false, $false, $false, $false, $true, 1, $false, "61-30=31", 2) | Out-Null
$d.Content.Find.Execute("59-59=0", $true, $false, $false, $false, $false, $true, 1, $false, "84-72=12", 2) | Out-Null
$d.Content.Find.Execute("81-35=46", $true, $false, $false, $false, $false, $true, 1, $false, "56+14=70", 2) | Out-Null
$d.Content.Find.Execute("64-25=39", $true, $false, $false, $false, $false, $true, 1, $false, "87-80=7", 2) | Out-Null
$d.Content.Find.Execute("25+44=69", $true, $false, $false, $false, $false, $true, 1, $false, "67-65=2", 2) | Out-Null
$d.Content.Find.Execute("50-33=17", $true, $false, $false, $false, $false, $true, 1, $false, "11+35=46", 2) | Out-Null
$d.Content.Find.Execute("19+67=86", $true, $false, $false, $false, $false, $true, 1, $false, "1+29=30", 2) | Out-Null
$d.Content.Find.Execute("39+3=42", $true, $false, $false, $false, $false, $true, 1, $false, "24-6=18", 2) | Out-Null
$d.Content.Find.Execute("61-14=47", $true, $false, $false, $false, $false, $true, 1, $false, "72-63=9", 2) | Out-Null
$d.Content.Find.Execute("74-20=54", $true, $false, $false, $false, $false, $true, 1, $false, "63-59=4", 2) | Out-Null
$d.Content.Find.Execute("91-8=83", $true, $false, $false, $false, $false, $true, 1, $false, "60-44=16", 2) | Out-Null
$d.Content.Find.Execute("56+10=66", $true, $false, $false, $false, $false, $true, 1, $false, "0+50=50", 2) | Out-Null
$d.Content.Find.Execute("23-1=22", $true, $false, $false, $false, $false, $true, 1, $false, "26+33=59", 2) | Out-Null
$d.Content.Find.Execute("15+38=53", $true, $false, $false, $false, $false, $true, 1, $false, "36+37=73", 2) | Out-Null
$d.Content.Find.Execute("89-40=49", $true, $false, $false, $false, $false, $true, 1, $false, "93-32=61", 2) | Out-Null
$d.Content.Find.Execute("41+36=77", $true, $false, $false, $false, $false, $true, 1, $false, "94-49=45", 2) | Out-Null
$d.Content.Find.Execute("24+36=60", $true, $false, $false, $false, $false, $true, 1, $false, "11+71=82", 2) | Out-Null
$d.Content.Find.Execute("90-22=68", $true, $false, $false, $false, $false, $true, 1, $false, "0+47=47", 2) | Out-Null
$d.Content.Find.Execute("50-47=3", $true, $false, $false, $false, $false, $true, 1, $false, "44+44=88", 2) | Out-Null
$d.Content.Find.Execute("4+7=11", $true, $false, $false, $false, $false, $true, 1, $false, "14+15=29", 2) | Out-Null
$d.Content.Find.Execute("13+27=40", $true, $false, $false, $false, $false, $true, 1, $false, "61-36=25", 2) | Out-Null
$d.Content.Find.Execute("47+4=51", $true, $false, $false, $false, $false, $true, 1, $false, "22+73=95", 2) | Out-Null
$d.Content.Find.Execute("43-5=38", $true, $false, $false, $false, $false, $true, 1, $false, "34+33=67", 2) | Out-Null
$d.Content.Find.Execute("97-73=24", $true, $false, $false, $false, $false, $true, 1, $false, "5+42=47", 2) | Out-Null
$d.Content.Find.Execute("64+18=82", $true, $false, $false, $false, $false, $true, 1, $false, "68-52=16", 2) | Out-Null
$d.Content.Find.Execute("57-12=45", $true, $false, $false, $false, $false, $true, 1, $false, "46-9=37", 2) | Out-Null
$d.Content.Find.Execute("27-24=3", $true, $false, $false, $false, $false, $true, 1, $false, "32+9=41", 2) | Out-Null
$d.Content.Find.Execute("66+9=75", $true, $false, $false, $false, $false, $true, 1, $false, "79-3=76", 2) | Out-Null
$d.Content.Find.Execute("77-57=20", $true, $false, $false, $false, $false, $true, 1, $false, "17-8=9", 2) | Out-Null
$d.Content.Find.Execute("80-62=18", $true, $false, $false, $false, $false, $true, 1, $false, "7+24=31", 2) | Out-Null
$d.Content.Find.Execute("26+43=69", $true, $false, $false, $false, $false, $true, 1, $false, "88-74=14", 2) | Out-Null
$d.Content.Find.Execute("7+17=24", $true, $false, $false, $false, $false, $true, 1, $false, "38-30=8", 2) | Out-Null
$d.Content.Find.Execute("96-92=4", $true, $false, $false, $false, $false, $true, 1, $false, "94-14=80", 2) | Out-Null
$d.Content.Find.Execute("79-23=56", $true, $false, $false, $false, $false, $true, 1, $false, "24+70=94", 2) | Out-Null
$d.Content.Find.Execute("78-39=39", $true, $false, $false, $false, $false, $true, 1, $false, "47-7=40", 2) | Out-Null
$d.Content.Find.Execute("48+49=97", $true, $false, $false, $false, $false, $true, 1, $false, "95-64=31", 2) | Out-Null
$d.Content.Find.Execute("31+42=73", $true, $false, $false, $false, $false, $true, 1, $false, "76+0=76", 2) | Out-Null
$d.Content.Find.Execute("26+20=46", $true, $false, $false, $false, $false, $true, 1, $false, "60-41=19", 2) | Out-Null
$d.Content.Find.Execute("1+82=83", $true, $false, $false, $false, $false, $true, 1, $false, "3+23=26", 2) | Out-Null
$d.Content.Find.Execute("90-34=56", $true, $false, $false, $false, $false, $true, 1, $false, "20+9=29", 2) | Out-Null
$d.Content.Find.Execute("0+31=31", $true, $false, $false, $false, $false, $true, 1, $false, "73-56=17", 2) | Out-Null
$d.Content.Find.Execute("90-18=72", $true, $false, $false, $false, $false, $true, 1, $false, "79-32=47", 2) | Out-Null
$d.Content.Find.Execute("16+79=95", $true, $false, $false, $false, $false, $true, 1, $false, "28+42=70", 2) | Out-Null
$d.Content.Find.Execute("16+38=54", $true, $false, $false, $false, $false, $true, 1, $false, "91-43=48", 2) | Out-Null
$d.Content.Find.Execute("55-54=1", $true, $false, $false, $false, $false, $true, 1, $false, "29-7=22", 2) | Out-Null
$d.Content.Find.Execute("20+72=92", $true, $false, $false, $false, $false, $true, 1, $false, "21+59=80", 2) | Out-Null
$d.Content.Find.Execute("82+4=86", $true, $false, $false, $false, $false, $true, 1, $false, "52-5=47", 2) | Out-Null
$d.Content.Find.Execute("87-18=69", $true, $false, $false, $false, $false, $true, 1, $false, "23+31=54", 2) | Out-Null
$d.Content.Find.Execute("89+1=90", $true, $false, $false, $false, $false, $true, 1, $false, "1+61=62", 2) | Out-Null
$d.Content.Find.Execute("29+65=94", $true, $false, $false, $false, $false, $true, 1, $false, "18+48=66", 2) | Out-Null
$d.Content.Find.Execute("50+44=94", $true, $false, $false, $false, $false, $true, 1, $false, "73-54=19", 2) | Out-Null
$d.Content.Find.Execute("97-17=80", $true, $false, $false, $false, $false, $true, 1, $false, "30+37=67", 2) | Out-Null
$d.Content.Find.Execute("90-72=18", $true, $false, $false, $false, $false, $true, 1, $false, "77-50=27", 2) | Out-Null
$d.Content.Find.Execute("94-22=72", $true, $false, $false, $false, $false, $true, 1, $false, "19+42=61", 2) | Out-Null
